$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6
$ws.Range("B3").Value = 11
$ws.Range("B4").Value = 26
$ws.Range("B5").Value = 3
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 1
$ws.Range("B9").Value = 14
$ws.Range("C9").Value = 2
$ws.Range("B10").Value = 22
$ws.Range("C10").Value = 3
$ws.Range("C11").Value = 3
$ws.Range("C12").Value = 1
$ws.Range("B24").Value = 20
$ws.Range("B25").Value = 56
$ws.Range("B26").Value = 23
$ws.Range("B27").Value = 50
$ws.Range("B30").Value = 17
$ws.Range("C30").Value = 2
$ws.Range("B31").Value = 43
$ws.Range("C31").Value = 6
$ws.Range("B32").Value = 9
$ws.Range("C32").Value = 3
$ws.Range("B33").Value = 14
$ws.Range("C33").Value = 4
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = 5
$ws.Range("B46").Value = 28
$ws.Range("B47").Value = 78
$ws.Range("B48").Value = 110
$ws.Range("B49").Value = 83
$ws.Range("B50").Value = 134
$ws.Range("B51").Value = 28
$ws.Range("B52").Value = 51
$ws.Range("C52").Value = 3
$ws.Range("B53").Value = 25
$ws.Range("C53").Value = 8
$ws.Range("B54").Value = 117
$ws.Range("C54").Value = 9
$ws.Range("B55").Value = 49
$ws.Range("C55").Value = 9
$ws.Range("B56").Value = 21
$ws.Range("C56").Value = 9
$ws.Range("B57").Value = 83
$ws.Range("B59").Value = 46
$ws.Range("B60").Value = 117
$ws.Range("B61").Value = 227
$ws.Range("B62").Value = 59
$ws.Range("B63").Value = 30
$ws.Range("C63").Value = 10
$ws.Range("B64").Value = 34
$ws.Range("C64").Value = 6
$ws.Range("B65").Value = 19
$ws.Range("C65").Value = 7
$ws.Range("B66").Value = 28
$ws.Range("C66").Value = 7
$ws.Range("B67").Value = 64
$ws.Range("C67").Value = 10
$ws.Range("B68").Value = 21
$ws.Range("B69").Value = 10
$ws.Range("B70").Value = 38
$ws.Range("B71").Value = 65
$ws.Range("B72").Value = 65
$ws.Range("B73").Value = 4
$ws.Range("B74").Value = 10
$ws.Range("C74").Value = 4
$ws.Range("C75").Value = 6
$ws.Range("B76").Value = 80
$ws.Range("C76").Value = 6
$ws.Range("B77").Value = 88
$ws.Range("B78").Value = 34
$ws.Range("C78").Value = 2
$ws.Range("B79").Value = 15
$ws.Range("B80").Value = 29
$ws.Range("B81").Value = 58
$ws.Range("B82").Value = 54
$ws.Range("B83").Value = 46
$ws.Range("B84").Value = 71
$ws.Range("B85").Value = 25
$ws.Range("C85").Value = 3
$ws.Range("B86").Value = 39
$ws.Range("C86").Value = 2
$ws.Range("B87").Value = 35
$ws.Range("C87").Value = 3
$ws.Range("B88").Value = 24
$ws.Range("C88").Value = 2
$ws.Range("B89").Value = 12
$ws.Range("C89").Value = 6